# Vehicle Ordering REST-Schnittstelle: rename service URLs from singular to
# plural (e.g. "services/order" -> "services/orders") per the updated
# Schnittstellendefinition, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "vehicle-ordering/services/types"
$ws.Range("C3").Value  = "vehicle-ordering/services/types/id"
$ws.Range("C4").Value  = "vehicle-ordering/services/models?type=id"
$ws.Range("C5").Value  = "vehicle-ordering/services/models/id"
$ws.Range("C6").Value  = "vehicle-ordering/services/customers?first=name&last=name"
$ws.Range("C7").Value  = "vehicle-ordering/services/customers/id"
$ws.Range("C8").Value  = "vehicle-ordering/services/orders?first=name&last=name&model=id"
$ws.Range("C9").Value  = "vehicle-ordering/services/orders/id"
$ws.Range("C10").Value = "vehicle-ordering/services/customers"
$ws.Range("C11").Value = "vehicle-ordering/services/orders"

# Move / leave the active selection on C10, matching the saved view state.
$ws.Range("C10").Select()
